# Apply odds updates to the active worksheet as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("S2").Value = 1.93
$ws.Range("T2").Value = 1.97
$ws.Range("AF2").Value = 34

# Row 3
$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 4.33
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 2.75
$ws.Range("U3").Value = 4.33
$ws.Range("V3").Value = 1.2
$ws.Range("AA3").Value = 6
$ws.Range("AH3").Value = 6
$ws.Range("AL3").Value = 10
$ws.Range("AO3").Value = 51

# Row 4
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("Q4").Value = 2.1
$ws.Range("R4").Value = 1.73

# Row 5
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 13
$ws.Range("Q5").Value = 1.8
$ws.Range("R5").Value = 2
$ws.Range("U5").Value = 3
$ws.Range("V5").Value = 1.36

# Row 6
$ws.Range("G6").Value = 2.5
$ws.Range("I6").Value = 2.8
$ws.Range("J6").Value = 3.1
$ws.Range("L6").Value = 3.4
$ws.Range("AD6").Value = 23
$ws.Range("AE6").Value = 19
$ws.Range("AM6").Value = 15
$ws.Range("AO6").Value = 29
